$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Append " (rejestracja przebiegów zrobiona)" to the paragraph that ends
#    with "...wcześniej zarejestrowanych przebiegów," and move the _GoBack
#    bookmark so that it sits right after the newly typed text (this mirrors
#    what Word does automatically after you type: the _GoBack bookmark marks
#    the last edited spot).
# ---------------------------------------------------------------------------

# Remove the existing _GoBack bookmark (it currently lives in the last,
# empty paragraph near the end of the document).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$target = $d.Content.Duplicate
$target.Find.Execute("wcześniej zarejestrowanych przebiegów,", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target.Collapse(0)
$target.InsertAfter(" (rejestracja przebiegów zrobiona)")

$insertEnd = $target.End

# Placing a zero-length bookmark exactly at "end of paragraph text" is
# mishandled when done directly, so pad with a throw-away character, anchor
# the bookmark next to it, then remove the padding -- the bookmark stays put.
$pad = $d.Range($insertEnd, $insertEnd)
$pad.InsertAfter("X")
$bmRange = $d.Range($insertEnd, $insertEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)
$d.Range($insertEnd, $insertEnd + 1).Delete()

# ---------------------------------------------------------------------------
# 2) Strike through the whole "(OKNO EDYCJI) mnożenie dwóch sygnałów," bullet
#    -- including its paragraph mark, so the list item's pPr/rPr also gets
#    the <w:strike/> toggle.
# ---------------------------------------------------------------------------

$hit = $d.Content.Duplicate
$hit.Find.Execute("(OKNO EDYCJI) mnożenie dwóch sygnałów,", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bulletPara = $hit.Paragraphs(1)
$bulletPara.Range.Font.StrikeThrough = 1
